{"js": "// Update the \"Version Control\" history table:\n//  - row \"2.8.1\": \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a \u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA) -> \u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM)\n//                 \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08      \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\n//  - row \"2.4.1\": \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a \u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\n//                 \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08      \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e23\u0e34\u0e28\u0e23\u0e32 (D)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Locate the two version-history rows by the version number in their first cell,\n// rather than hard-coding row indexes, so the script is resilient to any extra\n// rows elsewhere in the document.\nasync function findRowByFirstCellText(targetText) {\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n    const firstCell = cells.items[0];\n    firstCell.body.load(\"text\");\n    await context.sync();\n    if (firstCell.body.text.trim() === targetText) {\n      return cells;\n    }\n  }\n  return null;\n}\n\nasync function replaceInCell(cells, cellIndex, findText, replaceText) {\n  const cell = cells.items[cellIndex];\n  const results = cell.body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// Row for version 2.8.1\nconst row281Cells = await findRowByFirstCellText(\"2.8.1\");\nif (row281Cells) {\n  // \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a column (index 3)\n  await replaceInCell(row281Cells, 3, \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\", \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\");\n  await replaceInCell(row281Cells, 3, \" (QA)\", \" (DM)\");\n  // \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08 column (index 4)\n  await replaceInCell(row281Cells, 4, \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\");\n  await replaceInCell(row281Cells, 4, \"(SP)\", \" (TL)\");\n}\n\n// Row for version 2.4.1\nconst row241Cells = await findRowByFirstCellText(\"2.4.1\");\nif (row241Cells) {\n  // \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a column (index 3)\n  await replaceInCell(row241Cells, 3, \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\", \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\");\n  await replaceInCell(row241Cells, 3, \" (QA)\", \" (TL)\");\n  // \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08 column (index 4)\n  await replaceInCell(row241Cells, 4, \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", \"\u0e27\u0e23\u0e34\u0e28\u0e23\u0e32\");\n  await replaceInCell(row241Cells, 4, \"(SP)\", \" (D)\");\n}\n", "ps1": "# Update the \"Version Control\" history table:\n#  - row \"2.8.1\": \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a \u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA) -> \u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM)\n#                 \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08      \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\n#  - row \"2.4.1\": \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a \u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c (QA) -> \u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\n#                 \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08      \u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP) -> \u0e27\u0e23\u0e34\u0e28\u0e23\u0e32 (D)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Get-CellPlainText($cell) {\n    # Cell.Range.Text includes a trailing cell-mark (Chr 7) / paragraph marks;\n    # strip those so we can compare against plain version numbers.\n    return ($cell.Range.Text -replace \"[\\x07\\x0d]\", \"\")\n}\n\nfunction Replace-InCellRange($cell, [string]$findText, [string]$replaceText) {\n    $rng = $cell.Range\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceOne = 1\n    [void]$rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n\n# Locate the two version-history rows by the version number in their first cell,\n# rather than hard-coding row indexes, so the script is resilient to any extra\n# rows elsewhere in the document.\n$row281 = $null\n$row241 = $null\nfor ($i = 1; $i -le $tbl.Rows.Count; $i++) {\n    $cellText = Get-CellPlainText $tbl.Cell($i, 1)\n    if ($cellText -eq \"2.8.1\") { $row281 = $i }\n    if ($cellText -eq \"2.4.1\") { $row241 = $i }\n}\n\nif ($row281) {\n    # \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a column (index 4)\n    Replace-InCellRange $tbl.Cell($row281, 4) \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\" \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\"\n    Replace-InCellRange $tbl.Cell($row281, 4) \" (QA)\" \" (DM)\"\n    # \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08 column (index 5)\n    Replace-InCellRange $tbl.Cell($row281, 5) \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\n    Replace-InCellRange $tbl.Cell($row281, 5) \"(SP)\" \" (TL)\"\n}\n\nif ($row241) {\n    # \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a column (index 4)\n    Replace-InCellRange $tbl.Cell($row241, 4) \"\u0e13\u0e31\u0e10\u0e19\u0e31\u0e19\u0e17\u0e4c\" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\n    Replace-InCellRange $tbl.Cell($row241, 4) \" (QA)\" \" (TL)\"\n    # \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08 column (index 5)\n    Replace-InCellRange $tbl.Cell($row241, 5) \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" \"\u0e27\u0e23\u0e34\u0e28\u0e23\u0e32\"\n    Replace-InCellRange $tbl.Cell($row241, 5) \"(SP)\" \" (D)\"\n}\n"}
